$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-10-27 16:27:31"

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
